$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/12/2025  Through  5/18/2025"

# --- Crime Complaints table updates ---
# Row 14
$c = $ws.Range("M14")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -50
$c = $ws.Range("N14")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 0

# Row 15
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("L15")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 40

# Row 16
$c = $ws.Range("C16")
$c.NumberFormat = '#,##0'
$c.Value = 3
$c = $ws.Range("E16")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 0
$c = $ws.Range("F16")
$c.NumberFormat = '#,##0'
$c.Value = 11
$c = $ws.Range("G16")
$c.NumberFormat = '#,##0'
$c.Value = 11
$c = $ws.Range("H16")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 0
$c = $ws.Range("I16")
$c.NumberFormat = '#,##0'
$c.Value = 43
$c = $ws.Range("J16")
$c.NumberFormat = '#,##0'
$c.Value = 61
$c = $ws.Range("K16")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -29.508196721311
$c = $ws.Range("L16")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -36.764705882352
$c = $ws.Range("M16")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -23.214285714285
$c = $ws.Range("N16")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -88.154269972451

# Row 17
$c = $ws.Range("C17")
$c.NumberFormat = '#,##0'
$c.Value = 4
$c = $ws.Range("E17")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 33.333333333333
$c = $ws.Range("F17")
$c.NumberFormat = '#,##0'
$c.Value = 17
$c = $ws.Range("G17")
$c.NumberFormat = '#,##0'
$c.Value = 16
$c = $ws.Range("H17")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 6.25
$c = $ws.Range("I17")
$c.NumberFormat = '#,##0'
$c.Value = 93
$c = $ws.Range("J17")
$c.NumberFormat = '#,##0'
$c.Value = 93
$c = $ws.Range("K17")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 0
$c = $ws.Range("L17")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 8.139534883720
$c = $ws.Range("M17")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 144.736842105263
$c = $ws.Range("N17")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -21.186440677966

# Row 18
$c = $ws.Range("C18")
$c.NumberFormat = '#,##0'
$c.Value = 3
$c = $ws.Range("D18")
$c.NumberFormat = '#,##0'
$c.Value = 2
$c = $ws.Range("E18")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 50
$c = $ws.Range("F18")
$c.NumberFormat = '#,##0'
$c.Value = 13
$c = $ws.Range("G18")
$c.NumberFormat = '#,##0'
$c.Value = 12
$c = $ws.Range("H18")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 8.333333333333
$c = $ws.Range("I18")
$c.NumberFormat = '#,##0'
$c.Value = 69
$c = $ws.Range("J18")
$c.NumberFormat = '#,##0'
$c.Value = 57
$c = $ws.Range("K18")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 21.052631578947
$c = $ws.Range("L18")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -13.75
$c = $ws.Range("M18")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -19.767441860465
$c = $ws.Range("N18")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -92.023121387283

# Row 19
$c = $ws.Range("C19")
$c.NumberFormat = '#,##0'
$c.Value = 12
$c = $ws.Range("D19")
$c.NumberFormat = '#,##0'
$c.Value = 11
$c = $ws.Range("E19")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 9.090909090909
$c = $ws.Range("F19")
$c.NumberFormat = '#,##0'
$c.Value = 34
$c = $ws.Range("G19")
$c.NumberFormat = '#,##0'
$c.Value = 48
$c = $ws.Range("H19")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -29.166666666666
$c = $ws.Range("I19")
$c.NumberFormat = '#,##0'
$c.Value = 194
$c = $ws.Range("J19")
$c.NumberFormat = '#,##0'
$c.Value = 245
$c = $ws.Range("K19")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -20.816326530612
$c = $ws.Range("L19")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -23.320158102766
$c = $ws.Range("M19")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 38.571428571428
$c = $ws.Range("N19")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -65.232974910394

# Row 20
$c = $ws.Range("C20")
$c.NumberFormat = '#,##0'
$c.Value = 8
$c = $ws.Range("E20")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 33.333333333333
$c = $ws.Range("F20")
$c.NumberFormat = '#,##0'
$c.Value = 21
$c = $ws.Range("G20")
$c.NumberFormat = '#,##0'
$c.Value = 20
$c = $ws.Range("H20")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 5
$c = $ws.Range("I20")
$c.NumberFormat = '#,##0'
$c.Value = 93
$c = $ws.Range("J20")
$c.NumberFormat = '#,##0'
$c.Value = 90
$c = $ws.Range("K20")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 3.333333333333
$c = $ws.Range("L20")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 36.764705882352
$c = $ws.Range("M20")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 30.985915492957
$c = $ws.Range("N20")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -93.338108882521

# Row 21
$c = $ws.Range("C21")
$c.NumberFormat = '#,##0'
$c.Value = 30
$c = $ws.Range("E21")
$c.NumberFormat = '#,##0.00;"-"#,##0.00'
$c.Value = 20
$c = $ws.Range("F21")
$c.NumberFormat = '#,##0'
$c.Value = 97
$c = $ws.Range("H21")
$c.NumberFormat = '#,##0.00;"-"#,##0.00'
$c.Value = -9.345794392523
$c = $ws.Range("I21")
$c.NumberFormat = '#,##0'
$c.Value = 500
$c = $ws.Range("J21")
$c.NumberFormat = '#,##0'
$c.Value = 551
$c = $ws.Range("K21")
$c.NumberFormat = '#,##0.00;"-"#,##0.00'
$c.Value = -9.255898366606
$c = $ws.Range("L21")
$c.NumberFormat = '#,##0.00;"-"#,##0.00'
$c.Value = -11.190053285968
$c = $ws.Range("M21")
$c.NumberFormat = '#,##0.00;"-"#,##0.00'
$c.Value = 26.582278481012
$c = $ws.Range("N21")
$c.NumberFormat = '#,##0.00;"-"#,##0.00'
$c.Value = -84.880556395524

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = '#,##0'
$c.Value = 1
$c = $ws.Range("J22")
$c.NumberFormat = '#,##0'
$c.Value = 10
$c = $ws.Range("K22")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -90
$c = $ws.Range("L22")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -80
$c = $ws.Range("M22")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -66.666666666666

# Row 23
$c = $ws.Range("C23")
$c.NumberFormat = '#,##0'
$c.Value = 1
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "***.*"
$c = $ws.Range("F23")
$c.NumberFormat = '#,##0'
$c.Value = 4
$c = $ws.Range("H23")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 0
$c = $ws.Range("I23")
$c.NumberFormat = '#,##0'
$c.Value = 23
$c = $ws.Range("K23")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -20.689655172413
$c = $ws.Range("L23")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 0
$c = $ws.Range("M23")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 130

# Row 24
$c = $ws.Range("C24")
$c.NumberFormat = '#,##0'
$c.Value = 13
$c = $ws.Range("D24")
$c.NumberFormat = '#,##0'
$c.Value = 26
$c = $ws.Range("E24")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -50
$c = $ws.Range("F24")
$c.NumberFormat = '#,##0'
$c.Value = 73
$c = $ws.Range("G24")
$c.NumberFormat = '#,##0'
$c.Value = 82
$c = $ws.Range("H24")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -10.975609756097
$c = $ws.Range("I24")
$c.NumberFormat = '#,##0'
$c.Value = 417
$c = $ws.Range("J24")
$c.NumberFormat = '#,##0'
$c.Value = 405
$c = $ws.Range("K24")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 2.962962962962
$c = $ws.Range("L24")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -8.351648351648
$c = $ws.Range("M24")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 10.610079575596

# Row 25
$c = $ws.Range("C25")
$c.NumberFormat = '#,##0'
$c.Value = 4
$c = $ws.Range("D25")
$c.NumberFormat = '#,##0'
$c.Value = 8
$c = $ws.Range("E25")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -50
$c = $ws.Range("F25")
$c.NumberFormat = '#,##0'
$c.Value = 29
$c = $ws.Range("G25")
$c.NumberFormat = '#,##0'
$c.Value = 35
$c = $ws.Range("H25")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -17.142857142857
$c = $ws.Range("I25")
$c.NumberFormat = '#,##0'
$c.Value = 163
$c = $ws.Range("J25")
$c.NumberFormat = '#,##0'
$c.Value = 152
$c = $ws.Range("K25")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 7.236842105263
$c = $ws.Range("L25")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -23.831775700934

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = '#,##0'
$c.Value = 8
$c = $ws.Range("E26")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -12.5
$c = $ws.Range("F26")
$c.NumberFormat = '#,##0'
$c.Value = 35
$c = $ws.Range("G26")
$c.NumberFormat = '#,##0'
$c.Value = 42
$c = $ws.Range("H26")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -16.666666666666
$c = $ws.Range("I26")
$c.NumberFormat = '#,##0'
$c.Value = 185
$c = $ws.Range("J26")
$c.NumberFormat = '#,##0'
$c.Value = 166
$c = $ws.Range("K26")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 11.445783132530
$c = $ws.Range("L26")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 32.142857142857
$c = $ws.Range("M26")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 8.823529411764

# Row 27
$c = $ws.Range("C27")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("L27")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 16.666666666666

# Row 28
$c = $ws.Range("C28")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "***.*"
$c = $ws.Range("G28")
$c.NumberFormat = '#,##0'
$c.Value = 4
$c = $ws.Range("H28")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -25
$c = $ws.Range("L28")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 37.5

# Row 29
$c = $ws.Range("N29")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -94.736842105263

# Row 30
$c = $ws.Range("N30")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -93.333333333333

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "***.*"
$c = $ws.Range("F31")
$c.NumberFormat = '#,##0'
$c.Value = 1
$c = $ws.Range("G31")
$c.NumberFormat = '#,##0'
$c.Value = 2
$c = $ws.Range("H31")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -50
$c = $ws.Range("I31")
$c.NumberFormat = '#,##0'
$c.Value = 8
$c = $ws.Range("K31")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = -46.666666666666
$c = $ws.Range("L31")
$c.NumberFormat = '#,##0.0;"-"#,##0.0'
$c.Value = 33.333333333333
